# AutoDocGenR.xlsx edit
#
# 1. Clear cell F7 (previously held the shared string "autogen"), leaving its
#    existing cell style/formatting (s="12") intact but with no content.
# 2. Move the sheet's selection to F5 (was A11) and scroll the view so column
#    B is the first visible column (topLeftCell="B1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of F7 - keeps the cell's style, drops the t="s"/<v>
# shared-string reference so xl/sharedStrings.xml's usage count drops by one.
$ws.Range("F7").ClearContents()

# Scroll the window so column B is left-most in view, then select F5 to match
# the new activeCell/sqref of the saved sheetView.
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F5").Select()
